# Implemented the durablegoods endpoint call
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the date values for durable_goods, non_durable_goods, services rows,
# matching the style (number format) already used by B2 (monitored_goods row).
$ws.Range("B3").Value = 45125
$ws.Range("B4").Value = 45125
$ws.Range("B5").Value = 45125
$ws.Range("B3:B5").NumberFormat = $ws.Range("B2").NumberFormat

# Column sizing to match the bestFit widths captured in the diff (nearest
# achievable values given the engine's character-width rounding).
$ws.Columns.Item(1).ColumnWidth = 18
$ws.Columns.Item(2).ColumnWidth = 17.5

# Selection moved to E5 as recorded in the sheetView.
$ws.Range("E5").Select() | Out-Null
